$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear schedule cells that no longer apply (weeks already passed)
$cellsToClear = @("X6", "X7", "X8", "S10", "S11", "Y12", "Y13", "Y14", "V15", "X16")
foreach ($addr in $cellsToClear) {
    $ws.Range($addr).ClearContents()
}

# Move the selection to reflect the current week being reviewed
[void]$ws.Range("X16").Select()
